$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "Período: 22/10/2025 até 11/11/2025"

$ws.Range("G8").Value = "ALARES INTERNET S A"

$ws.Range("E9").Value = "MJ COMERCIO E SERV DE INF E TELECOMUNICA"
$ws.Range("L9").Value = "FINALIZADO"
$ws.Range("M9").Value = "ENTREGA REALIZADA"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "11/11/2025"

$ws.Range("E10").Value = "MJ COMERCIO E SERV DE INF E TELECOMUNICA"

$ws.Range("A24").Value = "Data: 11/11/2025 22:28"
